$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$h = $ws.Hyperlinks
$item = $h.Item(16)
Write-Host "pre count $($h.Count())"
[void]$item.Delete()
Write-Host "post count $($h.Count())"
$h2 = $ws.Hyperlinks
Write-Host "post count fresh $($h2.Count())"
